$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for 2022-Q4 at the top of
#    the data (row 2), push existing quarters down, and renumber the
#    running index in column A.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# copy column-A formatting (bold/bordered style) down onto the new row
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 15
$summary.Range("D2").Value = 1.6

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q4" sheet (positioned right after "总计",
#    before "2022-Q3") by duplicating "2022-Q2" - the closest existing
#    sheet in row-count/layout - then overwriting its cell values.
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($beforeSheet)
$q4 = $wb.Worksheets.Item("2022-Q2 (2)")
$q4.Name = "2022-Q4"

# the template only has 14 data rows (rows 2-15); Q4 needs 15 (rows 2-16)
# so clone the last data row's formatting down one more row first
$q4.Range("A15:H15").Copy($q4.Range("A16:H16"))

# columns that must stay text (fund codes keep leading zeros, the
# remaining metrics are stored as text in the source data too)
$q4.Range("B2:B16").NumberFormat = "@"
$q4.Range("D2:G16").NumberFormat = "@"

$q4Header = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$q4Data = @(
    @("002708", "大摩健康产业混合A", "23.18", "93.73", "3.61", "0.8368", "10"),
    @("200006", "长城消费增值混合", "5.59", "92.96", "4.89", "0.2734", "5"),
    @("014030", "大摩健康产业混合C", "4.20", "93.73", "3.61", "0.1516", "10"),
    @("000523", "国投瑞银医疗保健混合A", "2.11", "94.20", "3.75", "0.0791", "9"),
    @("014867", "摩根士丹利华鑫优悦安和混合C", "1.51", "92.74", "5.17", "0.0781", "10"),
    @("009893", "摩根士丹利华鑫优悦安和混合A", "1.17", "92.74", "5.17", "0.0605", "10"),
    @("001463", "光大保德信一带一路战略主题混合", "1.37", "85.11", "3.58", "0.0490", "4"),
    @("519097", "新华中小市值优选混合", "0.66", "70.51", "4.41", "0.0291", "2"),
    @("010703", "财通智选消费股票A", "0.45", "92.71", "3.02", "0.0136", "10"),
    @("010704", "财通智选消费股票C", "0.44", "92.71", "3.02", "0.0133", "10"),
    @("011082", "国投瑞银医疗保健混合C", "0.15", "94.20", "3.75", "0.0056", "9"),
    @("410009", "华富量子生命力混合", "0.10", "92.12", "4.62", "0.0046", "10"),
    @("001978", "泰信互联网+主题灵活配置混合", "0.06", "77.79", "2.05", "0.0012", "10"),
    @("006157", "财通量化核心优选混合", "0.09", "91.93", "1.08", "0.0010", "10"),
    @("519165", "新华鑫利灵活配置混合", "0.05", "22.54", "0.70", "0.0004", "10")
)

for ($c = 0; $c -lt $q4Header.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $q4Header[$c]
}

for ($r = 0; $r -lt $q4Data.Length; $r++) {
    $rowVals = $q4Data[$r]
    $q4.Cells.Item($r + 2, 1).Value = $r
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $q4.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}

# restore "总计" as the active sheet (unchanged from the original file)
$summary.Activate()
$null = $summary.Range("A1").Select()
